$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Orders")
$ws.Select()

$ws.Range("A9").Value = "Ipoh Coff"
$ws.Range("A10").Select()
